$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.078386306762695
$ws.Range("B1").Value = 1.170966625213623
$ws.Range("C1").Value = 0.920451819896698
$ws.Range("D1").Value = 5.111706733703613
$ws.Range("E1").Value = 1.987619876861572
